$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'58.010.42"
$ws.Range("E2").Value = '  -2.15%  '

# Row 3
$ws.Range("D3").Value = "'2.466.94"
$ws.Range("E3").Value = '  -2.37%  '

# Row 4
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").Value = "'518.11"
$ws.Range("E5").Value = '  -4.04%  '

# Row 6
$ws.Range("D6").Value = "'132.07"
$ws.Range("E6").Value = '  -4.41%  '

# Row 7
$ws.Range("E7").Value = '  +0.29%  '

# Row 8
$ws.Range("D8").Value = "'0.557"
$ws.Range("E8").Value = '  -1.88%  '

# Row 9
$ws.Range("D9").Value = "'0.0994"
$ws.Range("E9").Value = '  -2.35%  '

# Row 10
$ws.Range("E10").Value = '  -1.91%  '

# Row 11
$ws.Range("D11").Value = "'5.39"
$ws.Range("E11").Value = '  +0.61%  '

# Row 12
$ws.Range("D12").Value = "'0.341"
$ws.Range("E12").Value = '  -2.33%  '

# Row 13
$ws.Range("D13").Value = "'2.904.11"
$ws.Range("E13").Value = '  -1.49%  '

# Row 14
$ws.Range("D14").Value = "'57.914.24"
$ws.Range("E14").Value = '  -2.02%  '

# Row 15
$ws.Range("D15").Value = "'22.06"
$ws.Range("E15").Value = '  -5.05%  '

# Row 16
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("E16").Value = '  -2.47%  '

# Row 17
$ws.Range("D17").Value = "'2.475.18"

# Row 18
$ws.Range("D18").Value = "'10.86"
$ws.Range("E18").Value = '  -2.39%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = "'320.09"
$ws.Range("E19").Value = '  -1.93%  '

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = "'4.17"
$ws.Range("E20").Value = '  -3.26%  '

# Row 21
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = '  -0.08%  '

# Row 22
$ws.Range("E22").Value = '  -3.44%  '

# Row 23
$ws.Range("D23").Value = "'64.38"
$ws.Range("E23").Value = '  -1.45%  '

# Row 24
$ws.Range("D24").Value = "'0.409"
$ws.Range("E24").Value = '  -3.46%  '

# Row 25
$ws.Range("E25").Value = '  -0.43%  '

# Row 26
$ws.Range("E26").Value = '  -3.30%  '

# Row 27
$ws.Range("D27").Value = "'7.40"
$ws.Range("E27").Value = '  -3.73%  '

# Row 28
$ws.Range("D28").Value = "'0.0₃0747"

# Row 29
$ws.Range("D29").Value = "'6.38"
$ws.Range("E29").Value = '  -4.97%  '

# Row 30
$ws.Range("D30").Value = "'1.70"
$ws.Range("E30").Value = '  -4.83%  '

# Row 31
$ws.Range("D31").Value = "'164.91"
$ws.Range("E31").Value = '  -0.51%  '

# Row 32
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = '  -4.57%  '

# Row 33
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = '  +0.08%  '

# Row 34
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = '  -0.06%  '

# Row 35
$ws.Range("D35").Value = "'18.15"
$ws.Range("E35").Value = '  -1.97%  '

# Row 36
$ws.Range("D36").Value = "'1.33"
$ws.Range("E36").Value = '  -9.33%  '

# Row 37
$ws.Range("D37").Value = "'3.96"
$ws.Range("E37").Value = '  -4.35%  '

# Row 38
$ws.Range("E38").Value = '  -4.76%  '

# Row 39
$ws.Range("D39").Value = "'0.793"
$ws.Range("E39").Value = '  -3.54%  '

# Row 40
$ws.Range("E40").Value = '  -4.62%  '

# Row 41
$ws.Range("D41").Value = "'274.69"
$ws.Range("E41").Value = '  -5.12%  '

# Row 42
$ws.Range("D42").Value = "'4.97"
$ws.Range("E42").Value = '  -5.62%  '

# Row 43
$ws.Range("E43").Value = '  -3.39%  '

# Row 44
$ws.Range("D44").Value = "'127.01"
$ws.Range("E44").Value = '  -3.60%  '

# Row 45
$ws.Range("D45").Value = "'0.0910"
$ws.Range("E45").Value = '  -2.63%  '

# Row 46
$ws.Range("D46").Value = "'0.0490"
$ws.Range("E46").Value = '  -4.14%  '

# Row 47
$ws.Range("E47").Value = '  -3.44%  '

# Row 48
$ws.Range("D48").Value = "'17.07"
$ws.Range("E48").Value = '  -2.34%  '

# Row 49
$ws.Range("D49").Value = "'1.735.56"
$ws.Range("E49").Value = '  -1.74%  '

# Row 50
$ws.Range("D50").Value = "'0.971"
$ws.Range("E50").Value = '  -1.63%  '

# Row 51
$ws.Range("E51").Value = '  -2.70%  '
